$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 6728
$ws.Range("I3").Value = 7020
$ws.Range("H4").Value = 1675
$ws.Range("I4").Value = 1607
$ws.Range("I5").Value = 654
$ws.Range("I6").Value = 8148
$ws.Range("H7").Value = 25987
$ws.Range("I7").Value = 24157

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I3").Value = 64
$ws.Range("I7").Value = 279

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("I4").Value = 4
$ws.Range("I7").Value = 82

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I3").Value = 247
$ws.Range("I7").Value = 755

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I3").Value = 156
$ws.Range("I7").Value = 428

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 218
$ws.Range("I6").Value = 278
$ws.Range("I7").Value = 923

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I6").Value = 171
$ws.Range("I7").Value = 563

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 188
$ws.Range("I7").Value = 757
$ws.Range("I8").Value = 1445
$ws.Range("I10").Value = 178
$ws.Range("I11").Value = 368
$ws.Range("I15").Value = 279
$ws.Range("I16").Value = 70
$ws.Range("I19").Value = 679
$ws.Range("I23").Value = 235
$ws.Range("I29").Value = 1455
$ws.Range("I30").Value = 82
$ws.Range("I33").Value = 1072
$ws.Range("I34").Value = 110
$ws.Range("I36").Value = 331
$ws.Range("I37").Value = 755
$ws.Range("I41").Value = 105
$ws.Range("I42").Value = 889
$ws.Range("I44").Value = 182
$ws.Range("I47").Value = 175
$ws.Range("I48").Value = 308
$ws.Range("I49").Value = 161
$ws.Range("I50").Value = 126
$ws.Range("I51").Value = 287
$ws.Range("I52").Value = 549
$ws.Range("I53").Value = 264
$ws.Range("I54").Value = 483
$ws.Range("I55").Value = 279
$ws.Range("I61").Value = 24
$ws.Range("H63").Value = 226
$ws.Range("I63").Value = 79
$ws.Range("I65").Value = 563
$ws.Range("I66").Value = 68
$ws.Range("I67").Value = 923
$ws.Range("I70").Value = 40
$ws.Range("I71").Value = 71
$ws.Range("I78").Value = 326
$ws.Range("I79").Value = 686
$ws.Range("I83").Value = 519
$ws.Range("I85").Value = 1085
$ws.Range("I86").Value = 152
$ws.Range("I87").Value = 59
$ws.Range("I88").Value = 224
$ws.Range("I90").Value = 311
$ws.Range("I91").Value = 255
$ws.Range("I93").Value = 137
$ws.Range("I94").Value = 246
$ws.Range("I96").Value = 279
$ws.Range("I97").Value = 200
$ws.Range("I98").Value = 169
$ws.Range("I99").Value = 428
$ws.Range("H101").Value = 25987
$ws.Range("I101").Value = 24157

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I3").Value = 187
$ws.Range("I7").Value = 519

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 242
$ws.Range("I3").Value = 392
$ws.Range("I6").Value = 347
$ws.Range("I7").Value = 1072

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("I4").Value = 19
$ws.Range("I7").Value = 161

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I6").Value = 235
$ws.Range("I7").Value = 483

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 426
$ws.Range("I3").Value = 503
$ws.Range("I6").Value = 402
$ws.Range("I7").Value = 1455

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I3").Value = 199
$ws.Range("I6").Value = 216
$ws.Range("I7").Value = 679

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("I2").Value = 58
$ws.Range("I3").Value = 53
$ws.Range("I6").Value = 52
$ws.Range("I7").Value = 182

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I6").Value = 159
$ws.Range("I7").Value = 308

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 309
$ws.Range("I3").Value = 409
$ws.Range("I7").Value = 1085

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("I2").Value = 31
$ws.Range("I6").Value = 29
$ws.Range("I7").Value = 105

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I3").Value = 268
$ws.Range("I6").Value = 331
$ws.Range("I7").Value = 889

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("I6").Value = 82
$ws.Range("I7").Value = 178

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I2").Value = 73
$ws.Range("I4").Value = 46
$ws.Range("I7").Value = 326

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I6").Value = 85
$ws.Range("I7").Value = 279

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I5").Value = 7
$ws.Range("I7").Value = 235

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I2").Value = 79
$ws.Range("I4").Value = 15
$ws.Range("I7").Value = 255

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 199
$ws.Range("I7").Value = 686

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I2").Value = 95
$ws.Range("I3").Value = 110
$ws.Range("I7").Value = 331

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("I6").Value = 58
$ws.Range("I7").Value = 137

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I3").Value = 178
$ws.Range("I6").Value = 176
$ws.Range("I7").Value = 549

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("I3").Value = 31
$ws.Range("I4").Value = 12
$ws.Range("I7").Value = 110

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I6").Value = 141
$ws.Range("I7").Value = 246

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I3").Value = 53
$ws.Range("I7").Value = 175

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I6").Value = 106
$ws.Range("I7").Value = 279

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I6").Value = 107
$ws.Range("I7").Value = 169

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("I6").Value = 37
$ws.Range("I7").Value = 126

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("I6").Value = 29
$ws.Range("I7").Value = 68

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I2").Value = 146
$ws.Range("I4").Value = 35
$ws.Range("I7").Value = 368

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I4").Value = 17
$ws.Range("I7").Value = 188

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("I6").Value = 128
$ws.Range("I7").Value = 200

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("I3").Value = 12
$ws.Range("I7").Value = 40

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I2").Value = 66
$ws.Range("I3").Value = 76
$ws.Range("I7").Value = 224

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 428
$ws.Range("I3").Value = 419
$ws.Range("I5").Value = 45
$ws.Range("I7").Value = 1445

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("I4").Value = 72
$ws.Range("I7").Value = 152

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I6").Value = 110
$ws.Range("I7").Value = 311

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I2").Value = 62
$ws.Range("I7").Value = 287

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I6").Value = 125
$ws.Range("I7").Value = 264

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("I6").Value = 20
$ws.Range("I7").Value = 71

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I2").Value = 246
$ws.Range("I3").Value = 232
$ws.Range("I6").Value = 205
$ws.Range("I7").Value = 757

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("I6").Value = 34
$ws.Range("I7").Value = 59

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range("I6").Value = 10
$ws.Range("I7").Value = 24

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("I6").Value = 48
$ws.Range("I7").Value = 70

$wb.Save()